$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift the whole table 2 columns to the left (D:L -> B:J) ---
# Columns B:C are currently empty, deleting them shifts D:L into B:J.
$ws.Range("B:C").Delete()

# --- Step 2: insert a new blank column for "SBM Runtime" ---
# After the shift, the former "SSD Dissimilarity" column is now H; insert a
# fresh column there so the old H:J (SSD/MSQ/SSIM) move right to I:K.
$ws.Range("H:H").Insert()

# --- Step 3: rename header cells to match the new headings ---
$ws.Range("H5").Value = "SBM Runtime"
$ws.Range("I5").Value = "SSD Dissimilarity(Naïve-DP-SBM)"
$ws.Range("J5").Value = "MSQ (Naïve-DP-SBM)"
$ws.Range("K5").Value = "SSIM (Naïve-DP-SBM)"

# --- Step 4: column widths to match the new layout ---
$ws.Range("B1").ColumnWidth = 9.140625
$ws.Range("C1").ColumnWidth = 13.140625
$ws.Range("D1").ColumnWidth = 12.140625
$ws.Range("E1").ColumnWidth = 9
$ws.Range("F1").ColumnWidth = 19.140625
$ws.Range("G1").ColumnWidth = 13.42578125
$ws.Range("H1").ColumnWidth = 13.42578125
$ws.Range("I1").ColumnWidth = 17
$ws.Range("J1").ColumnWidth = 13.5703125
$ws.Range("K1").ColumnWidth = 14.42578125

# --- Step 5: view/selection state ---
$ws.Range("G9").Select()
